$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for all data rows (2-97)
# from 2023-10-25 (serial 45224) to 2023-11-03 (serial 45233).
$ws.Range("C2:C97").Value = 45233
